{"js": "// Word Online (Office.js) script \u2014 body of `async (context) => { ... }`.\n//\n// Changes applied (per the target diff):\n//  1. Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd) that sat\n//     right after the \"MP73010\" run in the title paragraph.\n//  2. Collapse the \">>  > your / stuff after this line >>>\" paragraph\n//     (previously split across 3 runs + proofErr markers) into a single\n//     run reading \">>>  your stuff after this line >>>\".\n//  3. Insert four new paragraphs right after that paragraph (and before\n//     the existing \"Ben changing things up!\" paragraph):\n//       - a blank paragraph\n//       - \"Version Control v1.00.1\"\n//       - \"Hello World!!\"\n//       - \"Sorry that\u2019s the best I could come up with.s\"\n\nconst body = context.document.body;\n\n// 1. Remove the _GoBack bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Find the paragraph that contains the \">>> ... >>>\" marker text and\n//    rewrite it as a single clean run (this also removes the proofErr\n//    tags that straddled the old run boundaries).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet markerParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"stuff after this line\") !== -1) {\n    markerParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!markerParagraph) {\n  throw new Error(\"Could not locate the '>>> ... stuff after this line >>>' paragraph\");\n}\n\nmarkerParagraph\n  .getRange()\n  .insertText(\">>>  your stuff after this line >>>\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Insert the four new paragraphs right after it, chaining each\n//    insertion off the previously inserted paragraph so they land in\n//    the correct order.\nconst blankParagraph = markerParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nconst versionParagraph = blankParagraph.insertParagraph(\n  \"Version Control v1.00.1\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst helloParagraph = versionParagraph.insertParagraph(\"Hello World!!\", Word.InsertLocation.after);\nawait context.sync();\n\nhelloParagraph.insertParagraph(\n  \"Sorry that\\u2019s the best I could come up with.s\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# Changes applied (per the target diff):\n#  1. Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd) that sat\n#     right after the \"MP73010\" run in the title paragraph.\n#  2. Collapse the \">>  > your / stuff after this line >>>\" paragraph\n#     (previously split across 3 runs + proofErr markers) into a single\n#     run reading \">>>  your stuff after this line >>>\".\n#  3. Insert four new paragraphs right after that paragraph (and before\n#     the existing \"Ben changing things up!\" paragraph):\n#       - a blank paragraph\n#       - \"Version Control v1.00.1\"\n#       - \"Hello World!!\"\n#       - \"Sorry that's the best I could come up with.s\"\n\n$d = $word.ActiveDocument\n\n# 1. Remove the _GoBack bookmark. (It's a hidden bookmark, so it won't\n#    show up in $d.Bookmarks.Count/enumeration, but direct access still\n#    works, same as real Word.)\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Rewrite the \">>> ... stuff after this line >>>\" paragraph as a\n#    single clean run. Using Find/Replace with the exact same text for\n#    both find & replacement collapses the three original runs (plus\n#    the proofErr gramStart/gramEnd markers) into one run.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \">>>  your stuff after this line >>>\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \">>>  your stuff after this line >>>\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 3. Find that paragraph again (by its text) and insert the four new\n#    paragraphs right after it, in order.\n$markerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*stuff after this line*\") {\n        $markerIndex = $i\n        break\n    }\n}\nif ($markerIndex -eq -1) {\n    throw \"Could not locate the '>>> ... stuff after this line >>>' paragraph\"\n}\n\n$markerPara = $d.Paragraphs($markerIndex)\n$markerPara.Range.InsertParagraphAfter()\n\n$blankPara = $d.Paragraphs($markerIndex + 1)\n$blankPara.Range.InsertParagraphAfter()\n\n$versionPara = $d.Paragraphs($markerIndex + 2)\n$versionPara.Range.InsertAfter(\"Version Control v1.00.1\")\n$versionPara.Range.InsertParagraphAfter()\n\n$helloPara = $d.Paragraphs($markerIndex + 3)\n$helloPara.Range.InsertAfter(\"Hello World!!\")\n$helloPara.Range.InsertParagraphAfter()\n\n$sorryPara = $d.Paragraphs($markerIndex + 4)\n$sorryPara.Range.InsertAfter(\"Sorry that\" + [char]0x2019 + \"s the best I could come up with.s\")\n"}
